$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 19

# Column A: date stored as plain text (matches the pattern used by the other
# recently-appended rows, e.g. A10:A18), so force text entry via a leading
# apostrophe, then strip the resulting "quote prefix" style so the cell ends
# up with no explicit style, just like its neighbours.
$cellA = $ws.Cells.Item($row, 1)
$cellA.Value = "'05/28/2025"
$cellA.Style = "Normal"

$ws.Cells.Item($row, 2).Value = 0.0004615799999999996
$ws.Cells.Item($row, 3).Value = 108323.5842107545
$ws.Cells.Item($row, 4).Value = 50
